$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 113, shifting existing rows 113:191 down to 114:192.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new weekly record.
$ws.Range("A113").Value2 = 11
$ws.Range("B113").Value2 = "Vega Monumental Concepción"
$ws.Range("C113").Value2 = "Bíobío"
$ws.Range("D113").Value2 = 44651
$ws.Range("D113").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E113").Value2 = 8
$ws.Range("F113").Value2 = 100112045
$ws.Range("G113").Value2 = "Zapallo"
$ws.Range("H113").Value2 = "Camote"
$ws.Range("I113").Value2 = "1a (cosecha)"
$ws.Range("J113").Value2 = 600
$ws.Range("K113").Value2 = 300
$ws.Range("L113").Value2 = 350
$ws.Range("M113").Value2 = 325
$ws.Range("N113").Value2 = "$/kilo (volumen en unidades)"
$ws.Range("O113").Value2 = "Región de O'Higgins"
$ws.Range("P113").Value2 = 325
$ws.Range("Q113").Value2 = 1
$ws.Range("R113").Value2 = "Hortaliza"
